# "Editing the Docs / Minor changes in Main_ID"
# The TO DO list gets a new leading "priority number" column, and the
# General_Info row is re-ordered (it now sits right before the Building_*
# rows and picks up the "V" flag), and Phaser_Output is reclassified from
# "Main" to "Subclass".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a brand-new column before column A - this pushes the existing
#    Main/Subclass/Building (old A), item name (old B) and V-flag (old C)
#    columns one slot to the right (new B, C, D).
$ws.Columns.Item(1).Insert()

# 1b. Re-apply the sheet's sort (same column, now shifted from C to D) so the
#     <sortState> bookkeeping follows the columns that moved - the row
#     contents are rewritten explicitly afterwards regardless of how the
#     engine reorders things here.
$sortRange = $ws.Range("B1:D14")
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("D1:D14")) | Out-Null
$sortObj.SetRange($sortRange)
$sortObj.Header = 0
$sortObj.Apply()

# 2. Final contents of the table (columns A:D), row by row, exactly as the
#    sheet should read after the edit.
$data = @(
    @(2,  "Subclass", "Directions",       "V"),
    @(4,  "Subclass", "Lanes",            "V"),
    @(1,  "Main",      "Main_Diagram",    "V"),
    @(3,  "Subclass", "Routes",           "V"),
    @(5,  "Subclass", "General_Info",     "V"),
    @(12, "Building",  "Building_Diagram", $null),
    @(13, "Building",  "Building_ID",      $null),
    @(14, "Building",  "Building_Table",   $null),
    @(6,  "Subclass", "LRT_Info",          $null),
    @(11, "Main",      "Main_ID",          $null),
    @(8,  "Main",      "Main_Table",       $null),
    @(7,  "Subclass", "Phaser_Output",     $null),
    @(10, "Subclass", "Phases",            $null),
    @(9,  "Subclass", "Section",           $null)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    if ($vals[3] -ne $null) {
        $ws.Cells.Item($row, 4).Value = $vals[3]
    } else {
        $ws.Cells.Item($row, 4).ClearContents()
    }
}

# 3. The conditional formatting used to watch columns A (class) and C
#    (V-flag); both slid one column right, so point the rules (still
#    anchored on the pre-insert B1:B14) at C1:C14, then fix up the
#    formulas to read from the new B/D columns. dxfId/priority/stopIfTrue
#    stay untouched so the existing style bindings keep working.
$fcs = $ws.Range("B1:B14").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("C1:C14"))
}

$fcs = $ws.Range("C1:C14").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $f = $fc.Formula1
    $f = $f.Replace("A1=", "B1=")
    $f = $f.Replace("C1=", "D1=")
    $fc.Formula1 = $f
}

# 4. Selection left where the edit session ended.
$ws.Range("A8").Select()

Write-Output "done"
